$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 odds updates (Jogos_da_Semana_FlashScore_2024-10-30)
$ws.Range("G3").Value = 2.45
$ws.Range("I3").Value = 2.75
$ws.Range("J3").Value = 3
$ws.Range("K3").Value = 2.3
$ws.Range("L3").Value = 3.2

$ws.Range("U3").Value = 1.53
$ws.Range("V3").Value = 2.38

$ws.Range("X3").Value = 15
$ws.Range("Y3").Value = 10
$ws.Range("Z3").Value = 26
$ws.Range("AA3").Value = 19
$ws.Range("AB3").Value = 23

$ws.Range("AD3").Value = 6.5

$ws.Range("AG3").Value = 126
$ws.Range("AH3").Value = 12
$ws.Range("AI3").Value = 15

$ws.Range("AL3").Value = 19

$ws.Range("AO3").Value = 13

$ws.Range("BB3").Value = 126
$ws.Range("BC3").Value = 401
